$d = $word.ActiveDocument

$d.Content.Find.Execute("40÷4=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷6=3, 1", 2)
$d.Content.Find.Execute("92÷7=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=16, 0", 2)
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2)
$d.Content.Find.Execute("82÷9=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "98÷6=16, 2", 2)
$d.Content.Find.Execute("57÷9=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "58÷6=9, 4", 2)
$d.Content.Find.Execute("11÷7=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "50÷7=7, 1", 2)
$d.Content.Find.Execute("22÷3=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "49÷6=8, 1", 2)
$d.Content.Find.Execute("53÷6=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "48÷6=8, 0", 2)
$d.Content.Find.Execute("83÷7=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "98÷6=16, 2", 2)
$d.Content.Find.Execute("69÷9=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "18÷7=2, 4", 2)
$d.Content.Find.Execute("52÷2=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2)
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "47÷6=7, 5", 2)
$d.Content.Find.Execute("67÷6=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=14, 3", 2)
$d.Content.Find.Execute("41÷5=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2)
$d.Content.Find.Execute("22÷8=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2)
$d.Content.Find.Execute("78÷8=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2)
$d.Content.Find.Execute("27÷4=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2)
$d.Content.Find.Execute("60÷6=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "23÷3=7, 2", 2)
$d.Content.Find.Execute("20÷8=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "65÷3=21, 2", 2)
$d.Content.Find.Execute("12÷6=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=13, 4", 2)
$d.Content.Find.Execute("68÷8=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2)
$d.Content.Find.Execute("30÷6=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=20, 1", 2)
$d.Content.Find.Execute("51÷6=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "17÷8=2, 1", 2)
$d.Content.Find.Execute("73÷4=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2)
$d.Content.Find.Execute("72÷3=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷7=6, 5", 2)
